# NIT-9018863266.xlsx — "Estado de Cuenta" update
#
# Commit: "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
#
# For this worker (CC 1143391503 - KELLY JOHANA PIEDRAHITA GRAJALES) a new
# overdue period (2509) is added to the statement, in addition to the
# existing one (2508). That means:
#   - "Cant. Periodos" (count of periods) goes from 1 to 2
#   - "VALOR MORA" (total overdue amount) doubles: 56940 -> 113880
#   - a new detail row is appended right below the existing 2508 row, with
#     the same worker/amounts but period 2509 - pushing the signature block
#     beneath it down by one row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the summary figures at the top of the statement -------------

# VALOR MORA (total overdue amount): 56940 -> 113880
$ws.Range("E11").Value = 113880

# Cant. Periodos (count of overdue periods): 1 -> 2
$ws.Range("F13").Value = 2

# --- Insert the new detail row for period 2509 ---------------------------

# Push row 21 (signature underline) and everything below it down by
# inserting a fresh row right after the existing detail row (16).
$ws.Rows("17:17").Insert()

# Clone the formatting of the existing detail row (16) onto the new row.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)

# Fill in the new row's values: same worker/amounts, new period (2509).
$ws.Range("B17").Value = $ws.Range("B16").Value2
$ws.Range("C17").Value = $ws.Range("C16").Value2
$ws.Range("D17").Value = $ws.Range("D16").Value2
$ws.Range("E17").Value = "2509"
$ws.Range("F17").Value = $ws.Range("F16").Value2
$ws.Range("G17").Value = $ws.Range("G16").Value2
$ws.Range("H17").Value = $ws.Range("H16").Value2
$ws.Range("I17").Value = $ws.Range("I16").Value2
$ws.Range("J17").Value = $ws.Range("J16").Value2

Write-Output "Updated VALOR MORA, Cant. Periodos, and inserted period 2509 row"
